# Updated cryptos list values (Price / Volume(1h)) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D) text, new Volume(1h) (column E) text.
# A value of $null means that column is unchanged for that row.
$updates = @(
    @{ Row = 2; D = '23.916.59'; E = '  -2.07%  ' },
    @{ Row = 3; D = '1.652.81'; E = '  -0.84%  ' },
    @{ Row = 4; D = '1.001'; E = '  +0.06%  ' },
    @{ Row = 5; D = '310.65'; E = '  -0.82%  ' },
    @{ Row = 6; D = $null; E = '  -0.02%  ' },
    @{ Row = 7; D = '0.3881'; E = '  -1.66%  ' },
    @{ Row = 8; D = '0.3817'; E = '  -2.46%  ' },
    @{ Row = 9; D = '51.80'; E = '  -0.69%  ' },
    @{ Row = 10; D = $null; E = '  -3.16%  ' },
    @{ Row = 11; D = '1.001'; E = '  +0.05%  ' },
    @{ Row = 12; D = '0.08466'; E = '  -1.30%  ' },
    @{ Row = 13; D = '24.00'; E = '  -1.48%  ' },
    @{ Row = 14; D = '7.084'; E = '  -2.88%  ' },
    @{ Row = 15; D = '8.087'; E = '  +1.78%  ' },
    @{ Row = 16; D = $null; E = '  -1.88%  ' },
    @{ Row = 17; D = '1.649.30'; E = '  -0.67%  ' },
    @{ Row = 18; D = '94.20'; E = '  -0.90%  ' },
    @{ Row = 19; D = $null; E = '  +0.05%  ' },
    @{ Row = 20; D = '19.68'; E = '  -4.22%  ' },
    @{ Row = 21; D = '6.969'; E = '  -0.23%  ' },
    @{ Row = 22; D = $null; E = '  +0.06%  ' },
    @{ Row = 23; D = $null; E = '  +0.41%  ' },
    @{ Row = 24; D = '23.906.34'; E = '  -2.05%  ' },
    @{ Row = 25; D = '2.438'; E = '  +0.46%  ' },
    @{ Row = 26; D = $null; E = '  -1.86%  ' },
    @{ Row = 27; D = '22.10'; E = '  -1.85%  ' },
    @{ Row = 28; D = '153.87'; E = '  -2.13%  ' },
    @{ Row = 29; D = '5.413'; E = '  -0.89%  ' },
    @{ Row = 30; D = '138.08'; E = '  -3.29%  ' },
    @{ Row = 31; D = '7.880'; E = '  -2.12%  ' },
    @{ Row = 32; D = '2.501'; E = '  -1.30%  ' },
    @{ Row = 33; D = '1.840.19'; E = '  -0.10%  ' },
    @{ Row = 34; D = $null; E = '  -3.45%  ' },
    @{ Row = 35; D = '0.08172'; E = '  -0.78%  ' },
    @{ Row = 36; D = '6.732'; E = '  -2.40%  ' },
    @{ Row = 37; D = $null; E = '  -3.52%  ' },
    @{ Row = 38; D = $null; E = '  -2.67%  ' },
    @{ Row = 39; D = '0.2682'; E = '  -2.93%  ' },
    @{ Row = 40; D = '0.09123'; E = '  -1.20%  ' },
    @{ Row = 41; D = '0.7585'; E = '  -1.65%  ' },
    @{ Row = 42; D = '13.54'; E = '  -1.89%  ' },
    @{ Row = 43; D = '1.426'; E = '  -1.60%  ' },
    @{ Row = 44; D = '16.49'; E = '  -0.13%  ' },
    @{ Row = 45; D = '0.6952'; E = '  -2.13%  ' },
    @{ Row = 46; D = '2.465'; E = '  -2.74%  ' },
    @{ Row = 47; D = '4.099'; E = '  -0.82%  ' },
    @{ Row = 48; D = '0.9999'; E = '  -0.05%  ' },
    @{ Row = 49; D = '0.08307'; E = $null },
    @{ Row = 50; D = '134.58'; E = '  -1.57%  ' },
    @{ Row = 51; D = '1.233'; E = '  -2.88%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the numeric-looking string as
        # text (matching the workbook's original inline-string cell type) instead
        # of silently converting it to a floating point number, then reset the
        # cell style so no stray quote-prefix formatting is left behind.
        $cell = $ws.Range("D" + $u.Row)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}

